# Update the battery charge worksheet with a new measurement row (row 4)
# and extend the four chart series to include the new data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Batterimåling")

# --- Copy formatting from row 3 into row 4 so the new cells pick up the
#     same number formats (date format for column A, 2-decimal format for
#     column E) without introducing brand new style/numFmt entries. ---
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("E3").Copy()
$ws.Range("E4").PasteSpecial(-4122)   # xlPasteFormats

# --- Add the new battery measurement values for 2025-01-04 ---
$ws.Range("A4").Value = 45661
$ws.Range("B4").Value = 12.31
$ws.Range("C4").Value = 11.8
$ws.Range("D4").Value = 12.84
$ws.Range("E4").Value = 13.12

# --- Extend each chart series so it covers the new row of data ---
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$valueCols = @("B", "C", "D", "E")
for ($i = 1; $i -le 4; $i++) {
    $series = $chart.SeriesCollection($i)
    $col = $valueCols[$i - 1]
    $series.Formula = "=SERIES(Batterimåling!`$$col`$1,Batterimåling!`$A`$2:`$A`$4,Batterimåling!`$$col`$2:`$$col`$4,$i)"
}

# --- Match the new active selection left behind by the edit ---
$ws.Range("E3").Select()
